# date time function changes
$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# SkillDetails sheet
# ----------------------------------------------------------------------
$wsSkill = $wb.Worksheets.Item("SkillDetails")

# Clear the Startdate value (H2) but keep its date format
$wsSkill.Range("H2").ClearContents()

# Enddate (I2): 10/28/2020 -> 12/28/2020
$wsSkill.Range("I2").Value = (Get-Date -Year 2020 -Month 12 -Day 28).Date

# StartTime on row 3 (J3): 7:00 PM -> 3:00 PM
$wsSkill.Range("J3").Value = 0.625

# Move the active selection to K5
$wsSkill.Range("K5").Select()

# ----------------------------------------------------------------------
# ManageListings sheet
# ----------------------------------------------------------------------
$wsListings = $wb.Worksheets.Item("ManageListings")

# Insert a new column before current column H (Startdate), for the new
# "Select day" field.
$wsListings.Range("H1").EntireColumn.Insert()

# New column header + value
$wsListings.Range("H2").Value = "Monday"
$wsListings.Range("H1").Value = "Select day"

# Category / SubCategory text changed
$wsListings.Range("C2").Value = "Writing & Translation"
$wsListings.Range("D2").Value = "Creative Writing"

# Startdate / Enddate (now shifted to I2/J2 after the column insert)
$wsListings.Range("I2").Value = (Get-Date -Year 2020 -Month 10 -Day 19).Date
$wsListings.Range("J2").Value = (Get-Date -Year 2020 -Month 10 -Day 25).Date

# Move the active selection to E11
$wsListings.Range("E11").Select()

$wb.Save()
